$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Hunk 1: insert two new empty paragraphs (styled "ListParagraph", no
# numbering) right after the "You should see the following page once
# you have verified your email:" paragraph, before the existing blank
# paragraph that follows it.
# ---------------------------------------------------------------------
$targetText = "You should see the following page once you have verified your email:"

$findSrc = $d.Content.Find
$findSrc.ClearFormatting()
$findSrc.Execute($targetText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$srcIndex = $findSrc.Parent.Paragraphs.First.Index

# Insert two paragraph breaks one at a time -- doing this in a single
# replacement with "^p^p" leaves a stray empty run behind on the
# in-between paragraph, whereas inserting them one at a time keeps both
# new paragraphs free of runs.
$findBreak1 = $d.Content.Find
$findBreak1.ClearFormatting()
$findBreak1.Execute($targetText, $false, $false, $false, $false, $false, $true, 1, $false, "$targetText^p", 2) | Out-Null

$findBreak2 = $d.Content.Find
$findBreak2.ClearFormatting()
$findBreak2.Execute($targetText, $false, $false, $false, $false, $false, $true, 1, $false, "$targetText^p", 2) | Out-Null

$newPara1 = $d.Paragraphs($srcIndex + 1)
$newPara2 = $d.Paragraphs($srcIndex + 2)

# Both paragraphs inherited the source paragraph's list numbering;
# strip it so only the ListParagraph style remains.
$newPara1.Range.ListFormat.RemoveNumbers()
$newPara2.Range.ListFormat.RemoveNumbers()
$newPara1.Style = "ListParagraph"
$newPara2.Style = "ListParagraph"

# ---------------------------------------------------------------------
# Hunk 2: trim the trailing sentence off the "(1) Change the URL..."
# paragraph, leaving just the closing period.
# ---------------------------------------------------------------------
$findTrim = $d.Content.Find
$findTrim.ClearFormatting()
$findTrim.Execute(" Or, alternatively, you can search for Digital Scholarship" + [char]8217 + "s account using the steps for option (2).", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------------
# Hunk 3 + 4: the "_GoBack" bookmark moves from the end of the document
# to the middle of the word "three" in "...rightmost button of three
# buttons...". Adding a bookmark with the same name automatically
# removes the old one (bookmark names are unique), satisfying both
# halves of the diff in a single step.
# ---------------------------------------------------------------------
$findSplit = $d.Content.Find
$findSplit.ClearFormatting()
$findSplit.Execute("rightmost button of t", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $findSplit.Parent.End
$bookmarkRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
